$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 79, pushing the existing rows 79-96 down to 80-97.
$ws.Rows("79:79").Insert()

# Populate the newly inserted row 79 with the new record.
$ws.Range("A79").Value = 7
$ws.Range("B79").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C79").Value = "Ñuble"
$ws.Range("D79").Value = 44642
$ws.Range("E79").Value = 16
$ws.Range("F79").Value = 100112030
$ws.Range("G79").Value = "Poroto granado"
$ws.Range("H79").Value = "Sin especificar"
$ws.Range("I79").Value = "Primera"
$ws.Range("J79").Value = 60
$ws.Range("K79").Value = 20000
$ws.Range("L79").Value = 21000
$ws.Range("M79").Value = 20500
$ws.Range("N79").Value = "$/saco 25 kilos"
$ws.Range("O79").Value = "Provincia de Diguillín"
$ws.Range("P79").Value = 820
$ws.Range("Q79").Value = 25
$ws.Range("R79").Value = "Hortaliza"
